# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '63.867.61'
$ws.Range('E2').Value = '  +8.17%  '
$ws.Range('D3').Value = '3.133.16'
$ws.Range('E3').Value = '  +5.78%  '
$ws.Range('E4').Value = '  -0.17%  '
Set-TextValue $ws.Range('D5') '588.28'
$ws.Range('E5').Value = '  +4.39%  '
Set-TextValue $ws.Range('D6') '146.16'
$ws.Range('E6').Value = '  +6.92%  '
$ws.Range('D8').Value = '3.124.66'
$ws.Range('E8').Value = '  +5.57%  '
Set-TextValue $ws.Range('D9') '0.534'
$ws.Range('E9').Value = '  +3.19%  '
Set-TextValue $ws.Range('D10') '0.153'
$ws.Range('E10').Value = '  +16.43%  '
$ws.Range('E11').Value = '  +9.09%  '
Set-TextValue $ws.Range('D12') '0.470'
$ws.Range('E12').Value = '  +5.14%  '
Set-TextValue $ws.Range('D13') '0.0000249'
$ws.Range('E13').Value = '  +9.67%  '
Set-TextValue $ws.Range('D14') '35.74'
$ws.Range('E14').Value = '  +6.89%  '
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = '3.636.63'
$ws.Range('E16').Value = '  +4.87%  '
Set-TextValue $ws.Range('D17') '7.18'
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').Value = '63.650.75'
$ws.Range('E18').Value = '  +7.49%  '
$ws.Range('D19').Value = '3.115.07'
$ws.Range('E19').Value = '  +4.91%  '
Set-TextValue $ws.Range('D20') '471.28'
$ws.Range('E20').Value = '  +8.69%  '
Set-TextValue $ws.Range('D21') '14.22'
$ws.Range('E21').Value = '  +5.21%  '
Set-TextValue $ws.Range('D22') '0.732'
$ws.Range('E22').Value = '  +2.20%  '
Set-TextValue $ws.Range('D23') '7.56'
$ws.Range('E23').Value = '  +8.46%  '
Set-TextValue $ws.Range('D24') '13.34'
$ws.Range('E24').Value = '  +2.23%  '
Set-TextValue $ws.Range('D25') '82.12'
$ws.Range('E25').Value = '  +3.23%  '
$ws.Range('E26').Value = '  +0.26%  '
Set-TextValue $ws.Range('D27') '8.72'
$ws.Range('E27').Value = '  +13.96%  '
$ws.Range('E28').Value = '  +1.04%  '
Set-TextValue $ws.Range('D29') '2.69'
$ws.Range('E29').Value = '  +5.90%  '
$ws.Range('E30').Value = '  -0.36%  '
Set-TextValue $ws.Range('D31') '6.91'
$ws.Range('E31').Value = '  +11.85%  '
Set-TextValue $ws.Range('D32') '27.15'
$ws.Range('E32').Value = '  +6.10%  '
$ws.Range('E33').Value = '  +6.45%  '
$ws.Range('D34').Value = '0.0₃0879'
$ws.Range('E34').Value = '  +15.68%  '
Set-TextValue $ws.Range('D35') '2.45'
$ws.Range('E35').Value = '  +19.04%  '
$ws.Range('E36').Value = '  +7.70%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D37') '3.36'
$ws.Range('E37').Value = '  +22.71%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D38') '6.12'
$ws.Range('E38').Value = '  +4.47%  '
Set-TextValue $ws.Range('D39') '50.84'
$ws.Range('E39').Value = '  +5.08%  '
Set-TextValue $ws.Range('D40') '445.80'
$ws.Range('E40').Value = '  +12.58%  '
Set-TextValue $ws.Range('D41') '8.73'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').Value = '2.932.45'
$ws.Range('E42').Value = '  +7.47%  '
Set-TextValue $ws.Range('D43') '0.0372'
$ws.Range('E43').Value = '  +6.57%  '
Set-TextValue $ws.Range('D44') '0.282'
$ws.Range('E44').Value = '  +13.44%  '
$ws.Range('E45').Value = '  +5.92%  '
Set-TextValue $ws.Range('D46') '2.19'
$ws.Range('E46').Value = '  +11.20%  '
Set-TextValue $ws.Range('D47') '35.93'
$ws.Range('E47').Value = '  +5.78%  '
Set-TextValue $ws.Range('D49') '123.73'
$ws.Range('E49').Value = '  +1.33%  '
Set-TextValue $ws.Range('D50') '0.112'
$ws.Range('E50').Value = '  +2.18%  '
Set-TextValue $ws.Range('D51') '24.83'
$ws.Range('E51').Value = '  +7.72%  '
